# Apply the commit's changes to the workbook via Excel COM interop.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Data change 1 -----------------------------------------------------
# The shared label used in column C for the second "ASC" group (rows
# 32-40, 42-50, 52-60) is renamed from "ASC2" to "ASC4". Updating these
# cells' values causes the dependent CONCAT formulas in column D
# (D32:D40, D42:D50, D52:D60) to recalculate their cached results
# (e.g. "ASC2_1" -> "ASC4_1") automatically.
$ws.Range("C32:C40").Value = "ASC4"
$ws.Range("C42:C50").Value = "ASC4"
$ws.Range("C52:C60").Value = "ASC4"

# --- Data change 2 -----------------------------------------------------
# Column G ("comparator" metric) for the 3rd sub-block of each of the
# two "ASC" groups (rows 23-30 and 53-60) is corrected from
# "earliestValue" to "penultimateValue" so it matches the first row of
# its own block (G22 / G52). Other blocks that also use "earliestValue"
# (rows 12-20, 42-50) are left untouched.
$ws.Range("G23:G30").Value = "penultimateValue"
$ws.Range("G53:G60").Value = "penultimateValue"

# --- View / selection change --------------------------------------------
# The sheet's window had scrolled to show row 10 at the top with cell
# M34 selected; it now shows row 34 at the top with G23:G30 selected
# (active cell G23).
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G23:G30").Select()
